$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. LM4040B10 voltage reference (row 20): correct LCSC part + description
$ws.Range("D20").Value = "C701887"
$ws.Range("E20").Value = "±0.5% 15mA Fixed SOT-23-3 Voltage References ROHS"

# Highlight the changed LCSC/Description cells for the LM4040 row with a
# thin outline (matches the "flag this row changed" styling used elsewhere
# in this BOM) and turn wrapping on for the pair.
$hl = $ws.Range("D20:E20")
$hl.WrapText = $true

$d20 = $ws.Range("D20")
$d20.Borders.Item(7).LineStyle = 1
$d20.Borders.Item(7).ColorIndex = 48
$d20.Borders.Item(8).LineStyle = 1
$d20.Borders.Item(8).ColorIndex = 48
$d20.Borders.Item(9).LineStyle = 1
$d20.Borders.Item(9).ColorIndex = 48
$d20.Borders.Item(10).LineStyle = -4142

$e20 = $ws.Range("E20")
$e20.Borders.Item(7).LineStyle = -4142
$e20.Borders.Item(8).LineStyle = 1
$e20.Borders.Item(8).ColorIndex = 48
$e20.Borders.Item(9).LineStyle = 1
$e20.Borders.Item(9).ColorIndex = 48
$e20.Borders.Item(10).LineStyle = 1
$e20.Borders.Item(10).ColorIndex = 47

# 2. New LED series resistors (R1,R5 / R0603 1.0k row, row 38): new LCSC part + description
$ws.Range("D38").Value = "C23228"
$ws.Range("E38").Value = "100mW Thick Film Resistors ±100ppm/℃ ±1% 680Ω 0603  Chip Resistor - Surface Mount ROHS"
